$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (alinea "1e"): progress crop from 0 to 100, now done and assigned
# to Bernardo (the original "-" shared formula in E12 is replaced by the
# literal name, same as typing it over the old formula cell).
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = "Bernardo"

# Reflect the author's final selection / scroll position in the sheet view.
$ws.Range("I18").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
